$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-02-14 Wednesday" "2024-02-15 Thursday"

Replace-Text "258×6=1548" "490×5=2450"
Replace-Text "247×2=494" "512×5=2560"
Replace-Text "363×8=2904" "155×6=930"
Replace-Text "974×9=8766" "864×5=4320"
Replace-Text "693×8=5544" "395×4=1580"
Replace-Text "587×3=1761" "422×6=2532"
Replace-Text "977×5=4885" "520×6=3120"
Replace-Text "599×8=4792" "432×6=2592"
Replace-Text "216×9=1944" "895×9=8055"
Replace-Text "902×6=5412" "825×6=4950"
Replace-Text "824×7=5768" "646×5=3230"
Replace-Text "460×9=4140" "151×5=755"
Replace-Text "833×3=2499" "888×2=1776"
Replace-Text "311×3=933" "938×4=3752"
Replace-Text "803×8=6424" "839×7=5873"
Replace-Text "912×9=8208" "189×8=1512"
Replace-Text "223×9=2007" "160×7=1120"
Replace-Text "987×3=2961" "650×9=5850"
Replace-Text "318×3=954" "683×7=4781"
Replace-Text "799×4=3196" "254×3=762"
Replace-Text "293×5=1465" "114×5=570"
Replace-Text "690×3=2070" "348×4=1392"
Replace-Text "748×6=4488" "587×4=2348"
Replace-Text "956×6=5736" "810×8=6480"
Replace-Text "350×3=1050" "752×6=4512"

Write-Output "Done"
